# Update default target info worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Data updates -------------------------------------------------------
# Row 13: Odroid M1 -> accelerator: opencl -> rknn, memory: 2G -> 4G
$ws.Range("F13").Value = "rknn"
$ws.Range("G13").Value = "4G"

# Row 14: Rasberry Pi5 -> accelerator: opencl -> tpu, memory: 2G -> 4G
$ws.Range("F14").Value = "tpu"
$ws.Range("G14").Value = "4G"

# Row 15: Comma 3X -> accelerator: opencl -> adreno
$ws.Range("F15").Value = "adreno"

# Row 16: KT cloud -> cpu: arm -> x86, accelerator: opencl -> cpu
$ws.Range("E16").Value = "x86"
$ws.Range("F16").Value = "cpu"

# --- Style updates --------------------------------------------------------
# Add a thin box border around A13:C18 (target_name / target info / engine columns)
$borderRange = $ws.Range("A13:C18")
$borderRange.Borders.LineStyle = 1
$borderRange.Borders.Weight = 2

# --- Column width -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30

# --- Selection ----------------------------------------------------------
$ws.Range("F30").Select()
